# Auto-generated edit script applying the scheduled-runner price refresh
# to the Leve profit tables across all sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 183511.64
$ws.Range("J17").Value = 183511.64
$ws.Range("L17").Value = 550534.92
$ws.Range("N17").Value = -550870.92
$ws.Range("H69").Value = 16992.215
$ws.Range("J69").Value = 16992.215
$ws.Range("L69").Value = 50976.645
$ws.Range("N69").Value = -52724.645
$ws.Range("H72").Value = 16992.215
$ws.Range("J72").Value = 16992.215
$ws.Range("L72").Value = 152929.935
$ws.Range("N72").Value = -161665.935
$ws.Range("H132").Value = 2279.276
$ws.Range("I132").Value = 2337.9167
$ws.Range("K132").Value = 7013.750100000001
$ws.Range("M132").Value = -4483.750100000001
$ws.Range("H136").Value = 76250
$ws.Range("J136").Value = 76250
$ws.Range("L136").Value = 76250
$ws.Range("N136").Value = -86450
$ws.Range("H138").Value = 2359.9138
$ws.Range("I138").Value = 1875.4857
$ws.Range("J138").Value = 3097.087
$ws.Range("K138").Value = 5626.4571
$ws.Range("L138").Value = 9291.261
$ws.Range("M138").Value = -486.4570999999996
$ws.Range("N138").Value = -19571.261
$ws.Range("H141").Value = 2947.4614
$ws.Range("I141").Value = 2943.0833
$ws.Range("K141").Value = 8829.249899999999
$ws.Range("M141").Value = -3649.249899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 465525.4
$ws.Range("I2").Value = 980928.5600000001
$ws.Range("J2").Value = 4375.2104
$ws.Range("K2").Value = 980928.5600000001
$ws.Range("L2").Value = 4375.2104
$ws.Range("M2").Value = -980815.5600000001
$ws.Range("N2").Value = -4601.2104
$ws.Range("H7").Value = 100000
$ws.Range("J7").Value = 100000
$ws.Range("L7").Value = 100000
$ws.Range("N7").Value = -100228
$ws.Range("H45").Value = 1621.7778
$ws.Range("I45").Value = 1537.125
$ws.Range("K45").Value = 1537.125
$ws.Range("M45").Value = -1160.125
$ws.Range("H61").Value = 33337198
$ws.Range("I61").Value = 35718140
$ws.Range("K61").Value = 35718140
$ws.Range("M61").Value = -35717928
$ws.Range("H74").Value = 30307692
$ws.Range("I74").Value = 43483750
$ws.Range("K74").Value = 43483750
$ws.Range("M74").Value = -43482876
$ws.Range("H77").Value = 30307692
$ws.Range("I77").Value = 43483750
$ws.Range("K77").Value = 217418750
$ws.Range("M77").Value = -217414382
$ws.Range("H116").Value = 465525.4
$ws.Range("I116").Value = 980928.5600000001
$ws.Range("J116").Value = 4375.2104
$ws.Range("K116").Value = 980928.5600000001
$ws.Range("L116").Value = 4375.2104
$ws.Range("M116").Value = -978634.5600000001
$ws.Range("N116").Value = -8963.2104
$ws.Range("H122").Value = 3281.875
$ws.Range("I122").Value = 2767.3333
$ws.Range("K122").Value = 8301.999899999999
$ws.Range("M122").Value = -5851.999899999999
$ws.Range("H132").Value = 3033812
$ws.Range("I132").Value = 3128599.8
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 9385799.399999999
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -9383269.399999999
$ws.Range("N132").Value = -6860
$ws.Range("H136").Value = 33337198
$ws.Range("I136").Value = 35718140
$ws.Range("K136").Value = 107154420
$ws.Range("M136").Value = -107151870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 465525.4
$ws.Range("I3").Value = 980928.5600000001
$ws.Range("J3").Value = 4375.2104
$ws.Range("K3").Value = 980928.5600000001
$ws.Range("L3").Value = 4375.2104
$ws.Range("M3").Value = -980814.5600000001
$ws.Range("N3").Value = -4603.2104
$ws.Range("H57").Value = 99999.25
$ws.Range("J57").Value = 99999.25
$ws.Range("L57").Value = 99999.25
$ws.Range("N57").Value = -101439.25
$ws.Range("H123").Value = 81999.2
$ws.Range("J123").Value = 96666
$ws.Range("L123").Value = 96666
$ws.Range("N123").Value = -106466
$ws.Range("H124").Value = 95608
$ws.Range("J124").Value = 95608
$ws.Range("L124").Value = 95608
$ws.Range("N124").Value = -105428
$ws.Range("H136").Value = 99999.25
$ws.Range("J136").Value = 99999.25
$ws.Range("L136").Value = 99999.25
$ws.Range("N136").Value = -110199.25
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 11826.777
$ws.Range("I22").Value = 25874.25
$ws.Range("K22").Value = 25874.25
$ws.Range("M22").Value = -25524.25
$ws.Range("H58").Value = 35724188
$ws.Range("I58").Value = 50011404
$ws.Range("K58").Value = 50011404
$ws.Range("M58").Value = -50011201
$ws.Range("H105").Value = 1539952.5
$ws.Range("I105").Value = 3334372
$ws.Range("K105").Value = 3334372
$ws.Range("M105").Value = -3332625
$ws.Range("H129").Value = 109999.25
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 109999.25
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 109999.25
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -119999.25
$ws.Range("H131").Value = 97464
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 97464
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 97464
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -107544
$ws.Range("H132").Value = 27027898
$ws.Range("I132").Value = 30303788
$ws.Range("K132").Value = 90911364
$ws.Range("M132").Value = -90908834
$ws.Range("H134").Value = 9260021
$ws.Range("I134").Value = 9616152
$ws.Range("K134").Value = 28848456
$ws.Range("M134").Value = -28845921
$ws.Range("H136").Value = 35724188
$ws.Range("I136").Value = 50011404
$ws.Range("K136").Value = 150034212
$ws.Range("M136").Value = -150031662
$ws.Range("H137").Value = 115198.8
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 111799
$ws.Range("J138").Value = 111799
$ws.Range("L138").Value = 111799
$ws.Range("N138").Value = -122079
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 168
$ws.Range("J46").Value = 100
$ws.Range("L46").Value = 300
$ws.Range("N46").Value = -482
$ws.Range("H121").Value = 44997.176
$ws.Range("I121").Value = 85226
$ws.Range("J121").Value = 1111.1818
$ws.Range("K121").Value = 255678
$ws.Range("L121").Value = 3333.5454
$ws.Range("M121").Value = -254368
$ws.Range("N121").Value = -5953.5454
$ws.Range("H137").Value = 4764405
$ws.Range("I137").Value = 6669312
$ws.Range("J137").Value = 2137.5
$ws.Range("K137").Value = 20007936
$ws.Range("L137").Value = 6412.5
$ws.Range("M137").Value = -20002836
$ws.Range("N137").Value = -16612.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11462.462
$ws.Range("I70").Value = 12099.167
$ws.Range("J70").Value = 10916.714
$ws.Range("K70").Value = 12099.167
$ws.Range("L70").Value = 10916.714
$ws.Range("M70").Value = -11829.167
$ws.Range("N70").Value = -11456.714
$ws.Range("H73").Value = 11462.462
$ws.Range("I73").Value = 12099.167
$ws.Range("J73").Value = 10916.714
$ws.Range("K73").Value = 12099.167
$ws.Range("L73").Value = 10916.714
$ws.Range("M73").Value = -11163.167
$ws.Range("N73").Value = -12788.714
$ws.Range("H107").Value = 649.12
$ws.Range("I107").Value = 441.45
$ws.Range("J107").Value = 1479.8
$ws.Range("K107").Value = 441.45
$ws.Range("L107").Value = 1479.8
$ws.Range("M107").Value = 1478.55
$ws.Range("N107").Value = -5319.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 475
$ws.Range("I93").Value = 475
$ws.Range("K93").Value = 475
$ws.Range("M93").Value = 773
$ws.Range("H95").Value = 30599.666
$ws.Range("I95").Value = 28000
$ws.Range("J95").Value = 31899.5
$ws.Range("K95").Value = 28000
$ws.Range("L95").Value = 31899.5
$ws.Range("M95").Value = -25254
$ws.Range("N95").Value = -37391.5
$ws.Range("H127").Value = 60715
$ws.Range("J127").Value = 60715
$ws.Range("L127").Value = 60715
$ws.Range("N127").Value = -70635
$ws.Range("H136").Value = 2813.8462
$ws.Range("I136").Value = 1530
$ws.Range("J136").Value = 3199
$ws.Range("K136").Value = 4590
$ws.Range("L136").Value = 9597
$ws.Range("M136").Value = -2040
$ws.Range("N136").Value = -14697

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I2").Value = 2750
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2750
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2638
$ws.Range("N2").ClearContents()
$ws.Range("H97").Value = 12098.833
$ws.Range("J97").Value = 12098.833
$ws.Range("L97").Value = 12098.833
$ws.Range("N97").Value = -14080.833
$ws.Range("H122").Value = 5787.5
$ws.Range("I122").Value = 5787.5
$ws.Range("K122").Value = 17362.5
$ws.Range("M122").Value = -14912.5
$ws.Range("H132").Value = 10207582
$ws.Range("I132").Value = 11631478
$ws.Range("K132").Value = 34894434
$ws.Range("M132").Value = -34891904
$ws.Range("H136").Value = 31251818
$ws.Range("I136").Value = 38462410
$ws.Range("K136").Value = 115387230
$ws.Range("M136").Value = -115384680
